# Weekly update: a new price record is added for the week of 2022-11-17,
# inserted as a new row 4 (pushing the existing rows 4-35 down to 5-36).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4; Excel shifts rows 4:35 down to 5:36
# and copies formatting (e.g. the date style on column D) from the row above.
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Range("A4").Value = 7
$ws.Range("B4").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C4").Value = "Ñuble"
$ws.Range("D4").Value = 44882
$ws.Range("E4").Value = 16
$ws.Range("F4").Value = 300000000
$ws.Range("G4").Value = "Espárragos"
$ws.Range("H4").Value = "Sin especificar"
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 1200
$ws.Range("K4").Value = 1000
$ws.Range("L4").Value = 1100
$ws.Range("M4").Value = 1050
$ws.Range("N4").Value = "`$/kilo"
$ws.Range("O4").Value = "Región de Ñuble"
$ws.Range("P4").Value = 1050
$ws.Range("Q4").Value = 1
$ws.Range("R4").Value = "Hortaliza"
